$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the existing "DropPack" header (column W) to the new column Y and
#    rename it to "DropPackList". Then put the new "MoveType" header in W
#    and the new "AtkDis" header in X.
$ws.Range("Y1").Value = "DropPackList"
$ws.Range("W1").Value = "MoveType"
$ws.Range("X1").Value = "AtkDis"

# Give the second part of the "AtkDis" header ("tkDis") its own font run,
# matching the rich-text formatting used in the source workbook.
$run = $ws.Range("X1").Characters(2, 5)
$run.Font.Name = "宋体"
$run.Font.Size = 11

# 2. Fill in the new MoveType (W) / AtkDis (X) data columns for the NPC rows.
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 20

$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 20

$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 20

$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 20

$ws.Range("W6").Value = 2
$ws.Range("X6").Value = 20

# 3. Adjust column widths: V shrinks from 32.125 to 25, the two new columns
#    (W, X) are 25 wide, and the new Y column is sized to fit its header.
$ws.Columns("V").ColumnWidth = 24.2857142857143
$ws.Columns("W").ColumnWidth = 24.2857142857143
$ws.Columns("X").ColumnWidth = 24.2857142857143
$ws.Columns("Y").ColumnWidth = 13.1428571428571

# 4. Update the view: scroll so column K is at the left edge and select X10.
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("X10").Select()
